# Commit: "Added homework for 04. MondoDB-and-Mongoose"
#
# 1) Remove the trailing "Homework" slide (sldId 286 / slide35.xml,
#    the last slide in the deck). Deleting it through the Slides
#    collection also removes its <p:sldId> entry from the
#    presentation's sldIdLst and the associated relationship /
#    content-type bookkeeping.
#
# 2) On the "Creating MongoDB client" code slide, the two adjacent
#    runs " " and "= new " (identical run formatting) get merged into
#    a single run " = new ".

$p = $ppt.ActivePresentation

# --- 1. Delete the "Homework" slide (sldId 286) -----------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 286) {
        $slide.Delete()
        break
    }
}

# --- 2. Merge the " " + "= new " runs on the MongoClient code slide --------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            $text = $tr.Text
            if ($text -eq "var mongoClient = new mongodb.MongoClient(server);") {
                # "var mongoClient = new mongodb.MongoClient(server);"
                #  123456789012345678...
                # The lone-space run sits right before "=", so characters
                # 16-22 (" = new ") span both former runs; rewriting that
                # sub-range collapses them into a single run while leaving
                # the surrounding runs untouched.
                $merged = $tr.Characters(16, 7)
                $merged.Text = " = new "
            }
        }
    }
}
